# "updated bug in date of afl" - the "Issued this " run had a stray
# trailing space baked into it (the very next run, the {{ current_date }}
# merge field, already supplies its own leading space), producing a
# double space in the rendered certificate ("Issued this  <date>").
# Trim that redundant trailing space from the "Issued this" run only -
# the following "{{ current_date }} " / "at NIA-Pangasinan ..." runs must
# stay untouched.

$d = $word.ActiveDocument

$text = $d.Content.Text
$idx = $text.IndexOf("Issued this")
if ($idx -lt 0) {
    throw "Could not find 'Issued this' in the document"
}

$dateMarker = "{{ current_date }} "
$cityMarker = "Urdaneta City."
$markerPos = $text.IndexOf($dateMarker, $idx)
$cityPos = $text.IndexOf($cityMarker, $idx)
if ($markerPos -lt 0 -or $cityPos -lt 0) {
    throw "Could not locate the surrounding runs for 'Issued this'"
}

# Remember the enclosing paragraph's end *before* editing, so the later
# run-resplit step can be scoped tightly to just this paragraph (which
# has no bold/italic runs of its own) instead of spilling into the rest
# of the document and disturbing unrelated formatting.
$cityEnd = $cityPos + $cityMarker.Length
$para = $d.Range($idx, $cityEnd).Paragraphs.Item(1)
$paraEnd = $para.Range.End

# 1) Rewrite the "Issued this " run's text without the trailing space.
#    (This engine re-tokenizes the paragraph's run list on any text
#    change, coalescing this run with the immediately-following
#    identically-formatted runs - corrected again in step 2 below.)
$runLen = $markerPos - $idx
$r1 = $d.Range($idx, $idx + $runLen)
$r1.Text = "Issued this"

# Re-measure: the document text shrank by exactly one character (the
# trimmed space), so every position from here on shifts left by one.
$shrink = $runLen - ("Issued this").Length
$afterRun1 = $idx + ("Issued this").Length
$afterRun2 = $markerPos - $shrink + $dateMarker.Length
$paraEndAdj = $paraEnd - $shrink

# 2) Re-split the coalesced run back into its original three runs by
#    toggling a formatting no-op (Bold on, then off) starting exactly at
#    each original run boundary. A pure formatting change does not
#    re-trigger the paragraph-wide run merge, so this carves the run
#    list back apart at those two points. Scoped to this paragraph only
#    (it contains no bold/italic text) so no other formatting is touched.
$rSplit1 = $d.Range($afterRun1, $paraEndAdj)
$rSplit1.Bold = 1
$rSplit1.Bold = 0

$rSplit2 = $d.Range($afterRun2, $paraEndAdj)
$rSplit2.Bold = 1
$rSplit2.Bold = 0
